# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F5, F6, F7, F9, F10, F12, F13
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value  = 581
$wsExhibit.Range("F6").Value  = 325
$wsExhibit.Range("F7").Value  = 2839
$wsExhibit.Range("F9").Value  = 8097
$wsExhibit.Range("F10").Value = 208
$wsExhibit.Range("F12").Value = 52
$wsExhibit.Range("F13").Value = 420

# Sheet "全部类型" (sheet4): F5, F6, F9, F11, F12, F14, F17
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 581
$wsAll.Range("F6").Value  = 325
$wsAll.Range("F9").Value  = 2839
$wsAll.Range("F11").Value = 8097
$wsAll.Range("F12").Value = 208
$wsAll.Range("F14").Value = 52
$wsAll.Range("F17").Value = 420
